# Generate Report for Handback
#
# The zh-cn and de-de handback files have now been processed and are in
# sync with en-US, so the localization-status report is regenerated:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#     (shown both on the per-language detail sheets and on the Overview sheet).
#   - The "Latest Handback DateTime" for each language is refreshed.
#   - The stale "handback file is not the latest" error is cleared.
#   - The Status / Error Detail columns are resized to fit the new content.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$dateFormat = "yyyy-mm-dd HH:mm:ss"
$hyperlinkColor = 15570276   # OLE (BGR) form of RGB FF6495ED
$xlUnderlineStyleSingle = 2

# ---------------------------------------------------------------------
# Overview sheet: zh-cn (col E) / de-de (col F) status cells
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.16
$overview.Columns.Item(6).ColumnWidth = 29.16

# restore formatting that the engine re-indexes whenever any cell changes
$overview.Range("B2").Font.Underline = $xlUnderlineStyleSingle
$overview.Range("B2").Font.Color = $hyperlinkColor
$overview.Range("G2").NumberFormat = $dateFormat

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-31 13:03:10"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.16
$zhcn.Columns.Item(16).ColumnWidth = 12.83

$zhcn.Range("A2").Font.Underline = $xlUnderlineStyleSingle
$zhcn.Range("A2").Font.Color = $hyperlinkColor
$zhcn.Range("I2").Font.Underline = $xlUnderlineStyleSingle
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("H2").NumberFormat = $dateFormat
$zhcn.Range("K2").NumberFormat = $dateFormat

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-31 13:03:21"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.16
$dede.Columns.Item(16).ColumnWidth = 12.83

$dede.Range("A2").Font.Underline = $xlUnderlineStyleSingle
$dede.Range("A2").Font.Color = $hyperlinkColor
$dede.Range("I2").Font.Underline = $xlUnderlineStyleSingle
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("H2").NumberFormat = $dateFormat
$dede.Range("K2").NumberFormat = $dateFormat
